$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in G1 with new "data up to" date
$ws.Range("G1").Value = "New sequences in the last 30 days (data up to2021-10-17)"

# Row 2 - Angola
$ws.Range("C2").Value = 447
$ws.Range("F2").Value = 79
$ws.Range("G2").Value = 1

# Row 3 - Botswana
$ws.Range("F3").Value = 67
$ws.Range("G3").Value = 1

# Row 4 - Democratic Republic of the Congo
$ws.Range("F4").Value = 118
$ws.Range("G4").Value = 1

# Row 5 - Eswatini
$ws.Range("F5").Value = 99
$ws.Range("G5").Value = 1

# Row 6 - Lesotho
$ws.Range("F6").Value = 271
$ws.Range("G6").Value = 1

# Row 7 - Madagascar
$ws.Range("F7").Value = 173
$ws.Range("G7").Value = 1

# Row 8 - Malawi
$ws.Range("F8").Value = 117
$ws.Range("G8").Value = 1

# Row 9 - Mauritius
$ws.Range("F9").Value = 121
$ws.Range("G9").Value = 1

# Row 10 - Mozambique
$ws.Range("F10").Value = 177
$ws.Range("G10").Value = 1

# Row 11 - Namibia
$ws.Range("C11").Value = 118
$ws.Range("F11").Value = 123
$ws.Range("G11").Value = 1

# Row 12 - South Africa
$ws.Range("C12").Value = 7116
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2021-09-25"
$ws.Range("F12").Value = 21
$ws.Range("G12").Value = 1

# Row 13 - Union of the Comoros
$ws.Range("F13").Value = 279
$ws.Range("G13").Value = 1

# Row 14 - Zambia
$ws.Range("F14").Value = 134
$ws.Range("G14").Value = 1

# Row 15 - Zimbabwe
$ws.Range("F15").Value = 113
$ws.Range("G15").Value = 1
